$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.948.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.791.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.55%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "622.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.788.46"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.59%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("E10").Value = "  +4.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.28"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.491"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000261"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.419.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.786.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.975.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.92%  "

$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("E19").Value = "  +1.49%  "

$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "509.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "

$ws.Range("E22").Value = "  +3.38%  "

$ws.Range("E23").Value = "  -2.33%  "

$ws.Range("E24").Value = "  +4.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000140"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +27.59%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  -1.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.59%  "

$ws.Range("E34").Value = "  -1.06%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.13%  "

$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.332"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.55%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.133"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.21%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "422.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.00%  "

$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.039.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0362"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.77%  "

$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.89%  "

$ws.Range("E51").Value = "  +0.98%  "
